# OPTIMIZACION DE DISEÑO NO.6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header C1 from "Número de DPI" to "DPI"
$ws.Range("C1").Value = "DPI"

# Add new header "Activo" in column J
$ws.Range("J1").Value = "Activo"

# Add boolean "Activo" values for the 4 data rows
$ws.Range("J2").Value = $true
$ws.Range("J3").Value = $true
$ws.Range("J4").Value = $false
$ws.Range("J5").Value = $true
